$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D2").Value = "'29.046.05"
$ws.Range("D2").Style = "Normal"

$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D3").Value = "'1.817.85"
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"

$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D5").Value = "'241.11"
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = "  -2.33%  "
$ws.Range("D6").Value = "'0.6135"
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("D8").Value = "'0.07314"
$ws.Range("D8").Style = "Normal"

$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D9").Value = "'0.2881"
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = "  -1.83%  "

$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D11").Value = "'0.07655"
$ws.Range("D11").Style = "Normal"

$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D12").Value = "'1.826.37"
$ws.Range("D12").Style = "Normal"

$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D14").Value = "'0.6582"
$ws.Range("D14").Style = "Normal"

$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D15").Value = "'81.40"
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = "  -3.94%  "
$ws.Range("D16").Value = "'0.000009001"
$ws.Range("D16").Style = "Normal"

$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("D17").Value = "'5.822"
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D18").Value = "'29.035.73"
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D19").Value = "'2.062.43"
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = "  +6.09%  "
$ws.Range("D20").Value = "'236.62"
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("D21").Value = "'12.41"
$ws.Range("D21").Style = "Normal"

$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D26").Value = "'0.1399"
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D27").Value = "'8.406"
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D28").Value = "'17.55"
$ws.Range("D28").Style = "Normal"

$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D29").Value = "'1.481"
$ws.Range("D29").Style = "Normal"

$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("D30").Value = "'0.05545"
$ws.Range("D30").Style = "Normal"

$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D31").Value = "'4.081"
$ws.Range("D31").Style = "Normal"

$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("D32").Value = "'4.086"
$ws.Range("D32").Style = "Normal"

$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D33").Value = "'1.207"
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D34").Value = "'0.7325"
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("D35").Value = "'1.811"
$ws.Range("D35").Style = "Normal"

$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D36").Value = "'1.127"
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D37").Value = "'2.616"
$ws.Range("D37").Style = "Normal"

$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("D38").Value = "'2.823"
$ws.Range("D38").Style = "Normal"

$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D39").Value = "'1.204.96"
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = "  -1.37%  "

$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("D41").Value = "'6.354"
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D42").Value = "'0.8904"
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D44").Value = "'100.72"
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D45").Value = "'1.968.18"
$ws.Range("D45").Style = "Normal"

$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D46").Value = "'64.29"
$ws.Range("D46").Style = "Normal"

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D47").Value = "'0.5083"
$ws.Range("D47").Style = "Normal"

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("E48").Value = "  -4.47%  "
$ws.Range("D48").Value = "'0.00000000118"
$ws.Range("D48").Style = "Normal"

$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D49").Value = "'0.3980"
$ws.Range("D49").Style = "Normal"

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D50").Value = "'8.985"
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = "  -1.21%  "
$ws.Range("D51").Value = "'0.05751"
$ws.Range("D51").Style = "Normal"

